# Updates the cryptocurrency price/volume table on Sheet1 to the refreshed
# snapshot values (GitHub Actions scheduled data refresh). Column D (Price)
# and E (Volume(1h)) are text-formatted strings, not numbers, so values that
# would otherwise be auto-parsed as numeric literals (losing significant
# trailing zeros, e.g. "100.00" -> 100) are written with a leading apostrophe
# to force Excel to keep them as text, exactly like typing them in the UI.
# Rows 36/37 (RenderToken / ARBITRUM) also swapped rank order in this refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.407.16'
$ws.Range("E2").Value = '  -1.87%  '
$ws.Range("D3").Value = '1.792.40'
$ws.Range("E3").Value = '  -2.31%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D6").Value = '''307.37'
$ws.Range("E6").Value = '  -1.04%  '
$ws.Range("D7").Value = '''0.4560'
$ws.Range("E7").Value = '  -1.35%  '
$ws.Range("D8").Value = '''0.3634'
$ws.Range("E8").Value = '  -0.72%  '
$ws.Range("D9").Value = '''47.24'
$ws.Range("E9").Value = '  +2.76%  '
$ws.Range("D10").Value = '''0.07085'
$ws.Range("E10").Value = '  -1.01%  '
$ws.Range("D11").Value = '''0.8759'
$ws.Range("E11").Value = '  -0.61%  '
$ws.Range("D12").Value = '''0.07884'
$ws.Range("E12").Value = '  +0.51%  '
$ws.Range("D13").Value = '''19.53'
$ws.Range("E13").Value = '  -0.49%  '
$ws.Range("D14").Value = '1.783.82'
$ws.Range("E14").Value = '  -2.66%  '
$ws.Range("D15").Value = '''5.273'
$ws.Range("E15").Value = '  -1.19%  '
$ws.Range("E16").Value = '  -0.79%  '
$ws.Range("D17").Value = '''85.03'
$ws.Range("E17").Value = '  -4.16%  '
$ws.Range("D18").Value = '''1.011'
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("D19").Value = '''0.000008533'
$ws.Range("E19").Value = '  -2.76%  '
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("D21").Value = '26.442.17'
$ws.Range("E21").Value = '  -1.85%  '
$ws.Range("D22").Value = '''14.24'
$ws.Range("E22").Value = '  -1.94%  '
$ws.Range("D23").Value = '''4.982'
$ws.Range("E23").Value = '  -0.50%  '
$ws.Range("D24").Value = '2.040.51'
$ws.Range("E24").Value = '  -1.64%  '
$ws.Range("D25").Value = '''10.50'
$ws.Range("E25").Value = '  +0.78%  '
$ws.Range("E26").Value = '  +0.77%  '
$ws.Range("E27").Value = '  +0.91%  '
$ws.Range("E28").Value = '  -1.74%  '
$ws.Range("D29").Value = '''2.039'
$ws.Range("E29").Value = '  +1.69%  '
$ws.Range("D30").Value = '''111.95'
$ws.Range("E30").Value = '  -1.52%  '
$ws.Range("D31").Value = '''4.853'
$ws.Range("E31").Value = '  -1.78%  '
$ws.Range("D32").Value = '''0.08670'
$ws.Range("E32").Value = '  -1.94%  '
$ws.Range("D33").Value = '''3.049'
$ws.Range("E33").Value = '  -1.86%  '
$ws.Range("D34").Value = '''4.447'
$ws.Range("E34").Value = '  -0.57%  '
$ws.Range("D35").Value = '''0.7267'
$ws.Range("E35").Value = '  -4.30%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '''1.109'
$ws.Range("E36").Value = '  -2.48%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '''2.643'
$ws.Range("E37").Value = '  -0.71%  '
$ws.Range("E38").Value = '  +0.30%  '
$ws.Range("D39").Value = '''1.076'
$ws.Range("E39").Value = '  -1.55%  '
$ws.Range("D40").Value = '''0.01941'
$ws.Range("E40").Value = '  +0.58%  '
$ws.Range("D41").Value = '''0.05109'
$ws.Range("E41").Value = '  -0.47%  '
$ws.Range("D42").Value = '''0.5271'
$ws.Range("E42").Value = '  +5.77%  '
$ws.Range("D43").Value = '''2.867'
$ws.Range("E43").Value = '  -2.23%  '
$ws.Range("D44").Value = '''6.897'
$ws.Range("E44").Value = '  -0.81%  '
$ws.Range("D45").Value = '''0.1517'
$ws.Range("E45").Value = '  -5.06%  '
$ws.Range("D46").Value = '''8.019'
$ws.Range("E46").Value = '  -4.20%  '
$ws.Range("D47").Value = '''0.4729'
$ws.Range("E47").Value = '  +1.21%  '
$ws.Range("E48").Value = '  +0.30%  '
$ws.Range("D49").Value = '''9.858'
$ws.Range("E49").Value = '  -3.73%  '
$ws.Range("D50").Value = '''100.00'
$ws.Range("E50").Value = '  -2.57%  '
$ws.Range("D51").Value = '''1.586'
$ws.Range("E51").Value = '  -1.73%  '
